$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 228, shifting existing rows 228:335 down to 229:336
$ws.Rows(228).Insert()

# Fill in the new row 228 with its data
$ws.Cells.Item(228, 1).Value = 10
$ws.Cells.Item(228, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(228, 3).Value = "La Araucanía"
$ws.Cells.Item(228, 4).Value = 44466
$ws.Cells.Item(228, 5).Value = 9
$ws.Cells.Item(228, 6).Value = "Fruta"
$ws.Cells.Item(228, 7).Value = 100108
$ws.Cells.Item(228, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(228, 9).Value = 100108006
$ws.Cells.Item(228, 10).Value = "Plátano"
$ws.Cells.Item(228, 11).Value = "Sin especificar"
$ws.Cells.Item(228, 12).Value = "Pintón"
$ws.Cells.Item(228, 13).Value = 1750
$ws.Cells.Item(228, 14).Value = 16000
$ws.Cells.Item(228, 15).Value = 18000
$ws.Cells.Item(228, 16).Value = 16686
$ws.Cells.Item(228, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(228, 18).Value = "Ecuador"
$ws.Cells.Item(228, 19).Value = 834
$ws.Cells.Item(228, 20).Value = 20
